$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D17").Value = 44526
$ws.Range("M17").Value = 250
$ws.Range("P17").Value = 3200
$ws.Range("S17").Value = 1600
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 3000
$ws.Range("O18").Value = 3500
$ws.Range("P18").Value = 3250
$ws.Range("S18").Value = 1625
$ws.Range("D19").Value = 44218
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 2500
$ws.Range("O19").Value = 2500
$ws.Range("P19").Value = 2500
$ws.Range("S19").Value = 1250
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 4500
$ws.Range("O20").Value = 5000
$ws.Range("P20").Value = 4750
$ws.Range("S20").Value = 2375
$ws.Range("D21").Value = 44516
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 4000
$ws.Range("O21").Value = 4000
$ws.Range("P21").Value = 4000
$ws.Range("R21").Value = 'Región de Ñuble'
$ws.Range("S21").Value = 2000
$ws.Range("D22").Value = 44168
$ws.Range("N22").Value = 5000
$ws.Range("O22").Value = 5500
$ws.Range("P22").Value = 5250
$ws.Range("R22").Value = 'Provincia de Curicó'
$ws.Range("S22").Value = 2625
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 3000
$ws.Range("O23").Value = 3500
$ws.Range("P23").Value = 3250
$ws.Range("S23").Value = 1625
$ws.Range("D24").Value = 44203
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 100
$ws.Range("N24").Value = 2500
$ws.Range("O24").Value = 2500
$ws.Range("P24").Value = 2500
$ws.Range("S24").Value = 1250
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 3000
$ws.Range("O25").Value = 3500
$ws.Range("P25").Value = 3250
$ws.Range("S25").Value = 1625
$ws.Range("D26").Value = 44217
$ws.Range("L26").Value = 'Segunda'
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 2500
$ws.Range("P26").Value = 2500
$ws.Range("R26").Value = 'Región de Ñuble'
$ws.Range("S26").Value = 1250
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 200
$ws.Range("N27").Value = 2000
$ws.Range("O27").Value = 2500
$ws.Range("P27").Value = 2250
$ws.Range("S27").Value = 1125
$ws.Range("D28").Value = 44202
$ws.Range("L28").Value = 'Segunda'
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = 1500
$ws.Range("O28").Value = 1500
$ws.Range("P28").Value = 1500
$ws.Range("S28").Value = 750
$ws.Range("D29").Value = 44167
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = 3500
$ws.Range("O29").Value = 4000
$ws.Range("P29").Value = 3750
$ws.Range("R29").Value = 'Región de O''Higgins'
$ws.Range("S29").Value = 1875
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 3000
$ws.Range("O30").Value = 3500
$ws.Range("P30").Value = 3250
$ws.Range("S30").Value = 1625
$ws.Range("D31").Value = 44208
$ws.Range("L31").Value = 'Segunda'
$ws.Range("M31").Value = 50
$ws.Range("N31").Value = 2500
$ws.Range("O31").Value = 2500
$ws.Range("P31").Value = 2500
$ws.Range("S31").Value = 1250
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 200
$ws.Range("N32").Value = 3000
$ws.Range("O32").Value = 3500
$ws.Range("P32").Value = 3250
$ws.Range("S32").Value = 1625
$ws.Range("D33").Value = 44211
$ws.Range("N33").Value = 2500
$ws.Range("O33").Value = 2500
$ws.Range("P33").Value = 2500
$ws.Range("R33").Value = 'Región de Ñuble'
$ws.Range("S33").Value = 1250
$ws.Range("D34").Value = 44505
$ws.Range("L34").Value = 'Segunda'
$ws.Range("M34").Value = 100
$ws.Range("N34").Value = 7000
$ws.Range("O34").Value = 7000
$ws.Range("P34").Value = 7000
$ws.Range("R34").Value = 'Provincia de Curicó'
$ws.Range("S34").Value = 3500
$ws.Range("L35").Value = 'Primera'
$ws.Range("M35").Value = 200
$ws.Range("N35").Value = 3000
$ws.Range("O35").Value = 3500
$ws.Range("P35").Value = 3250
$ws.Range("S35").Value = 1625
$ws.Range("D36").Value = 44204
$ws.Range("L36").Value = 'Segunda'
$ws.Range("M36").Value = 100
$ws.Range("N36").Value = 2500
$ws.Range("O36").Value = 2500
$ws.Range("P36").Value = 2500
$ws.Range("R36").Value = 'Región de Ñuble'
$ws.Range("S36").Value = 1250
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 400
$ws.Range("N37").Value = 3000
$ws.Range("O37").Value = 3500
$ws.Range("P37").Value = 3250
$ws.Range("S37").Value = 1625
$ws.Range("D38").Value = 44189
$ws.Range("L38").Value = 'Segunda'
$ws.Range("M38").Value = 200
$ws.Range("N38").Value = 2500
$ws.Range("O38").Value = 2500
$ws.Range("P38").Value = 2500
$ws.Range("R38").Value = 'Provincia de Curicó'
$ws.Range("S38").Value = 1250
$ws.Range("D39").Value = 44209
$ws.Range("M39").Value = 50
$ws.Range("N39").Value = 3000
$ws.Range("O39").Value = 3000
$ws.Range("P39").Value = 3000
$ws.Range("S39").Value = 1500
